$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value (45188 -> 45189) for every
# data row from row 2 through row 342. Update the whole range in one go.
$ws.Range("C2:C342").Value = 45189
